$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 41, shifting existing rows 41-85 down to 42-86.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new record.
$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value = "Arica y Parinacota"
$ws.Range("D41").Value = 45272
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = 100112052
$ws.Range("G41").Value = "Albahaca"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 350
$ws.Range("K41").Value = 800
$ws.Range("L41").Value = 1000
$ws.Range("M41").Value = 886
$ws.Range("N41").Value = "$/paquete"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 886
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"
